# Primary Member Details screen TC added
# Adds two new worksheets (AddPrimaryMember, AddDependent) with sample test
# data laid out as Excel Tables, mirroring the other "Add*" test-data sheets
# already in the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the two new worksheets at the end of the workbook (after "AddPM")
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPrimary = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsPrimary.Name = "AddPrimaryMember"

$wsDependent = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsPrimary)
$wsDependent.Name = "AddDependent"

# ---------------------------------------------------------------------------
# 2. AddPrimaryMember - headers + one sample row, turned into "Table2"
# ---------------------------------------------------------------------------
$wsPrimary.Range("A1").Value = "UserName"
$wsPrimary.Range("B1").Value = "Password"
$wsPrimary.Range("C1").Value = "Cust_Num"
$wsPrimary.Range("D1").Value = "Firstname"
$wsPrimary.Range("E1").Value = "Lastname"
$wsPrimary.Range("F1").Value = "emailaddress"
$wsPrimary.Range("G1").Value = "P_Address"
$wsPrimary.Range("H1").Value = "ZiPCode"
$wsPrimary.Range("I1").Value = "Agee"
$wsPrimary.Range("J1").Value = "Gender"
$wsPrimary.Range("K1").Value = "Maritalstatus"
$wsPrimary.Range("L1").Value = "SSN"

$wsPrimary.Range("A2").Value = "nmusallam"
$wsPrimary.Range("B2").Value = "don1thack"
$wsPrimary.Range("C2").Value = "53569687"
$wsPrimary.Range("D2").Value = "soha"
$wsPrimary.Range("E2").Value = "ashraf"
$wsPrimary.Range("F2").Value = "t_7ryr@test.com"
$wsPrimary.Range("G2").Value = "18 test street"
$wsPrimary.Range("H2").Value = 75202
$wsPrimary.Range("I2").Value = 25
$wsPrimary.Range("J2").Value = "female"
$wsPrimary.Range("K2").Value = "single"
$wsPrimary.Range("L2").Value = "9"

$loPrimary = $wsPrimary.ListObjects.Add(1, $wsPrimary.Range("A1:L2"), 0, 1)
$loPrimary.Name = "Table2"
$loPrimary.TableStyle = "TableStyleLight15"

$wsPrimary.PageSetup.Orientation = 1
$wsPrimary.Activate()
$wsPrimary.Range("M7").Select()

# ---------------------------------------------------------------------------
# 3. AddDependent - headers + four sample rows, turned into "Table1"
# ---------------------------------------------------------------------------
$wsDependent.Range("A1").Value = "UserName"
$wsDependent.Range("B1").Value = "Password"
$wsDependent.Range("C1").Value = "CustomerNum"
$wsDependent.Range("D1").Value = "PrimaryNum"
$wsDependent.Range("E1").Value = "FirstName"
$wsDependent.Range("F1").Value = "LastName"
$wsDependent.Range("G1").Value = "Relationship"
$wsDependent.Range("H1").Value = "Age"
$wsDependent.Range("I1").Value = "Gender"
$wsDependent.Range("J1").Value = "Eligible"
$wsDependent.Range("K1").Value = "Eligible Reason"
$wsDependent.Range("L1").Value = "ProofDate"

$wsDependent.Range("A2").Value = "nmusallam"
$wsDependent.Range("B2").Value = "don1thack"
$wsDependent.Range("C2").Value = 53570720
$wsDependent.Range("D2").Value = 105698752
$wsDependent.Range("E2").Value = "DependentTest1"
$wsDependent.Range("F2").Value = "Testyo"
$wsDependent.Range("G2").Value = "child"
$wsDependent.Range("H2").Value = 30
$wsDependent.Range("I2").Value = "male"
$wsDependent.Range("J2").Value = "No"
$wsDependent.Range("K2").Value = "Disability"
$wsDependent.Range("L2").Value = 44197
$wsDependent.Range("L2").NumberFormat = "m/d/yyyy"

$wsDependent.Range("A3").Value = "nmusallam"
$wsDependent.Range("B3").Value = "don1thack"
$wsDependent.Range("C3").Value = 53570720
$wsDependent.Range("D3").Value = 105698752
$wsDependent.Range("E3").Value = "DependentTest2"
$wsDependent.Range("F3").Value = "Testyo"
$wsDependent.Range("G3").Value = "child"
$wsDependent.Range("H3").Value = 15
$wsDependent.Range("I3").Value = "female"
$wsDependent.Range("J3").Value = "No"
$wsDependent.Range("K3").Value = "Under Age"
$wsDependent.Range("L3").Value = "NA"

$wsDependent.Range("A4").Value = "nmusallam"
$wsDependent.Range("B4").Value = "don1thack"
$wsDependent.Range("C4").Value = 53570720
$wsDependent.Range("D4").Value = 105698752
$wsDependent.Range("E4").Value = "DependentTest3"
$wsDependent.Range("F4").Value = "Testyo"
$wsDependent.Range("G4").Value = "spouse"
$wsDependent.Range("H4").Value = 20
$wsDependent.Range("I4").Value = "male"
$wsDependent.Range("J4").Value = "NA"
$wsDependent.Range("K4").Value = "NA"
$wsDependent.Range("L4").Value = "NA"

$wsDependent.Range("A5").Value = "nmusallam"
$wsDependent.Range("B5").Value = "don1thack"
$wsDependent.Range("C5").Value = 53570720
$wsDependent.Range("D5").Value = 105698752
$wsDependent.Range("E5").Value = "DependentTest3"
$wsDependent.Range("F5").Value = "Testyo"
$wsDependent.Range("G5").Value = "child"
$wsDependent.Range("H5").Value = 19
$wsDependent.Range("I5").Value = "male"
$wsDependent.Range("J5").Value = "No"
$wsDependent.Range("K5").Value = "Married"
$wsDependent.Range("L5").Value = "NA"

$loDependent = $wsDependent.ListObjects.Add(1, $wsDependent.Range("A1:L5"), 0, 1)
$loDependent.Name = "Table1"
$loDependent.TableStyle = "TableStyleLight1"

$wsDependent.Range("I13").Select()

# ---------------------------------------------------------------------------
# 4. Minor selection changes on pre-existing sheets captured by the diff
# ---------------------------------------------------------------------------
$wsMandatory = $wb.Worksheets.Item("Mandatory")
$wsMandatory.Activate()
$wsMandatory.Range("D23").Select()

$wsAddPM = $wb.Worksheets.Item("AddPM")
$wsAddPM.Activate()
$wsAddPM.Range("J2").Select()

# ---------------------------------------------------------------------------
# 5. Leave "AddPrimaryMember" as the final active sheet/tab
# ---------------------------------------------------------------------------
$wsPrimary.Activate()
